# Insert a new weekly price record as row 291 in the "Alcachofa" sheet,
# pushing the previously-existing rows 291-321 down to 292-322 (the last
# of which, former row 321, becomes the new row 322).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 291..321 down one position, growing the used range to A1:R322.
$ws.Rows.Item(291).Insert()

# Populate the newly inserted row 291 with the new observation.
$ws.Cells.Item(291, 1).Value  = 10
$ws.Cells.Item(291, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(291, 3).Value  = "La Araucanía"
$ws.Cells.Item(291, 4).Value  = 45132
$ws.Cells.Item(291, 5).Value  = 9
$ws.Cells.Item(291, 6).Value  = 100112013
$ws.Cells.Item(291, 7).Value  = "Alcachofa"
$ws.Cells.Item(291, 8).Value  = "Madrigal"
$ws.Cells.Item(291, 9).Value  = "Primera"
$ws.Cells.Item(291, 10).Value = 55
$ws.Cells.Item(291, 11).Value = 14000
$ws.Cells.Item(291, 12).Value = 14000
$ws.Cells.Item(291, 13).Value = 14000
$ws.Cells.Item(291, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(291, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(291, 16).Value = 350
$ws.Cells.Item(291, 17).Value = 40
$ws.Cells.Item(291, 18).Value = "Hortaliza"
